# Commit: "user can see all sanity with url"
#
# 1. Fix a typo in the "End Serial" value of row 7 (sheet "Sheet1"):
#    AB0000000000000000000000001100 -> AC0000000000000000000000001100
#
# 2. Append a new data row (row 8) to Sheet1 with a 7th sanity-check entry
#    (Row=7, Reference Number=45646, Description="row7",
#     Start Serial="AA00000000000000000000000000101",
#     End Serial="AA0000000000000000000000000105", Date=7/5/2012)
#
# 3. Leave the active selection on the newly added row (E8).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. correct the existing row 7 "End Serial" value ---
$ws1.Range("E7").Value = "AC0000000000000000000000001100"

# --- 2. add the new row 8 ---
$ws1.Range("A8").Value = 7
$ws1.Range("B8").Value = 45646
$ws1.Range("C8").Value = "row7"
$ws1.Range("D8").Value = "AA00000000000000000000000000101"
$ws1.Range("E8").Value = "AA0000000000000000000000000105"
$ws1.Range("F8").Value = 41095
$ws1.Range("F8").NumberFormat = "MM/DD/YY"

# --- 3. move the selection to the new row, like the author left it ---
$ws1.Range("E8").Select() | Out-Null
